# Re-pull data: update column F (dSF) values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    5  = 10
    7  = -1
    8  = 2
    10 = -1
    11 = -1
    12 = 2
    13 = 6
    14 = 1
    15 = 1
    16 = 1
    17 = -4
    19 = 3
    21 = 2
    22 = 1
    23 = -1
    24 = -1
    25 = 2
    26 = 3
    27 = -3
    29 = 1
    31 = -1
    32 = 2
    33 = -3
    34 = -1
    35 = -2
    36 = 1
    37 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
